{"js": "// Insert three new bullet paragraphs right after the \"Research & Data\n// Analytics Leadership\" paragraph under the Siege Analytics / PARTNER\n// role, matching the diff: new bullets about voter-file discovery,\n// boundary estimation algorithm, and cost savings.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the exact paragraph that reads \"Research & Data Analytics Leadership\".\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Research & Data Analytics Leadership\") {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not find paragraph \"Research & Data Analytics Leadership\"');\n}\n\n// New bullet lines to insert, in order, immediately after the anchor paragraph.\nconst newBullets = [\n  \"\\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n  \"\\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n  \"\\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations $5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n];\n\n// Insert each bullet right after the anchor, one by one, so they end up\n// in the same order as in the diff (each new insertion becomes the new\n// \"after\" point for the next one).\nlet insertAfter = anchor;\nfor (const text of newBullets) {\n  insertAfter = insertAfter.insertParagraph(text, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Insert three new bullet paragraphs right after the \"Research & Data\n# Analytics Leadership\" paragraph under the Siege Analytics / PARTNER\n# role, matching the diff: new bullets about voter-file discovery,\n# boundary estimation algorithm, and cost savings.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`n\", [char]7)\n    if ($t -eq \"Research & Data Analytics Leadership\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find paragraph 'Research & Data Analytics Leadership'\"\n}\n\n$newBullets = @(\n    \"\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n    \"\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n    \"\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n)\n\n# Walk a range forward from the end of the anchor paragraph, inserting a\n# new paragraph mark plus its text for each bullet, then stepping over the\n# freshly-inserted paragraph mark so the next bullet lands after it (this\n# keeps the bullets in the same order as the diff).\n$r = $target.Range\n$r.Collapse(0)  # wdCollapseEnd\n\nforeach ($bullet in $newBullets) {\n    $r.InsertParagraphAfter()\n    $r.Collapse(0)\n    $r.InsertAfter($bullet)\n    $r.Collapse(0)\n    $r.Move(1, 1) | Out-Null  # wdCharacter: step past the new paragraph mark\n}\n"}
